$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "We are assuming the information suggested by the EUA model"
